$d = $word.ActiveDocument

# --- Title (paragraph 1): rephrase from a statement to a question and
#     swap "Irrigation ... Fertilization ... Rice" for "Technologies ... Rates ... Affect Rice".
$title = $d.Paragraphs(1).Range
$title.Find.Execute(
    "Effects of Alternate Wetting and Drying Irrigation and Nitrogen Fertilization on Sheath Blight of Rice",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Do Alternate Wetting and Drying Irrigation Technologies and Nitrogen Rates Affect Rice Sheath Blight?",
    2)

# --- Author line (paragraph 2): add Ole (B.O. Sander) as a third author.
$authors = $d.Paragraphs(2).Range
$authors.Find.Execute(
    "N.P. Castilla",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "N.P. Castilla and B.O. Sander",
    2)

# --- Affiliation note (paragraph 3): "second author:" -> "second and third authors:"
$affil = $d.Paragraphs(3).Range
$affil.Find.Execute(
    "second author:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "second and third authors:",
    2)
